$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.084.55"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").Value = "3.611.73"
$ws.Range("E3").Value = "  +3.17%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "195.62"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.47%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.207"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.46%  "
$ws.Range("E10").Value = "  -0.41%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.97"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.83%  "
$ws.Range("E12").Value = "  +0.39%  "
$ws.Range("E13").Value = "  -0.30%  "
$ws.Range("D14").Value = "4.177.31"
$ws.Range("E14").Value = "  +2.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "13.16"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.10%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "593.88"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.31"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.70%  "
$ws.Range("D18").Value = "70.293.64"
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("D19").Value = "3.607.13"
$ws.Range("E19").Value = "  +3.27%  "
$ws.Range("E20").Value = "  +1.69%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.997"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.93"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.52%  "
$ws.Range("E23").Value = "  +2.71%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "102.90"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.18%  "
$ws.Range("E25").Value = "  +0.43%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.08"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.88"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.58%  "
$ws.Range("E28").Value = "  -2.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "34.15"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.76%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.13"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.23%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.31"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.32"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.83%  "
$ws.Range("E33").Value = "  +0.70%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.60"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.19%  "
$ws.Range("D35").Value = "3.911.61"
$ws.Range("E35").Value = "  +5.97%  "
$ws.Range("D36").Value = "0.0₃0849"
$ws.Range("E36").Value = "  +5.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.18"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.73%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "529.10"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.79%  "
$ws.Range("E39").Value = "  +0.12%  "
$ws.Range("B40").Value = "InjectiveProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.21"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.15%  "
$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.394"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.35%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.57"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.52%  "
$ws.Range("E43").Value = "  -2.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0458"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.86"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.27%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.141"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.84%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.34"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.13%  "
$ws.Range("E48").Value = "  -1.85%  "
$ws.Range("E49").Value = "  -0.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000249"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.88%  "
$ws.Range("E51").Value = "  +3.36%  "
